$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# Overview sheet: Status mirrored into the per-locale columns (zh-cn / de-de)
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

# zh-cn sheet: Status column, new handback datetime, cleared error detail
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus
$zhcn.Range("K2").Value = "2016-10-21 03:47:35"
$zhcn.Range("K3").Value = "2016-10-21 03:47:35"
$zhcn.Range("P2").Value = ""
$zhcn.Range("P3").Value = ""

# de-de sheet: Status column, new handback datetime, cleared error detail
$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus
$dede.Range("K2").Value = "2016-10-21 03:47:53"
$dede.Range("K3").Value = "2016-10-21 03:47:53"
$dede.Range("P2").Value = ""
$dede.Range("P3").Value = ""

# Column width adjustments (Status / Error Detail columns widened / narrowed
# to fit the new report text) across all three sheets
$overview.Columns.Item(5).ColumnWidth = 29.9777050018311
$overview.Columns.Item(6).ColumnWidth = 29.9777050018311

$zhcn.Columns.Item(3).ColumnWidth = 29.9777050018311
$zhcn.Columns.Item(16).ColumnWidth = 13.7470531463623

$dede.Columns.Item(3).ColumnWidth = 29.9777050018311
$dede.Columns.Item(16).ColumnWidth = 13.7470531463623
